# Adds the new match row (row 47) to the bottom of the Moldova Super Liga
# 2023-2024 results sheet, mirroring the style of the row above it (row 46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from the last existing data row (46) down to
# the new row (47) so the new cells pick up the same styles (bold/border for
# column A, the custom date number format for column E, plain for the rest)
# instead of minting brand-new style records.
$ws.Range("A46:V46").Copy()
$ws.Range("A47:V47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row's values.
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "moldova"
$ws.Range("C47").Value = "super-liga"
$ws.Range("D47").Value = "2023-2024"
$ws.Range("E47").Value = 45241.5
$ws.Range("F47").Value = "Sparta Selemet"
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = "Balti"
$ws.Range("I47").Value = 5
$ws.Range("J47").Value = 6.41
$ws.Range("K47").Value = "11/11/2023 00:43"
$ws.Range("L47").Value = 8.109999999999999
$ws.Range("M47").Value = "11/11/2023 11:58"
$ws.Range("N47").Value = 4.53
$ws.Range("O47").Value = "11/11/2023 00:43"
$ws.Range("P47").Value = 4.79
$ws.Range("Q47").Value = "11/11/2023 11:58"
$ws.Range("R47").Value = 1.34
$ws.Range("S47").Value = "11/11/2023 00:43"
$ws.Range("T47").Value = 1.29
$ws.Range("U47").Value = "11/11/2023 11:58"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/moldova/super-liga/sparta-selemet-csf-balti/0xzfnnB2/"
